$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Range("A2").Value = "Última actualización: 04:53:50"
$ws1.Range("A3").Value = "Total filas: 32"

# Insert two new rows before the current row 17, shifting existing data down
$ws1.Range("A17:A18").EntireRow.Insert()

# Fill the two newly inserted rows
$ws1.Range("A17").Value = "04:53:50"
$ws1.Range("B17").Value = "04:54"
$ws1.Range("C17").Value = "11_ETCHEVERRY"
$ws1.Range("D17").Value = 1
$ws1.Range("E17").Value = "LP1912"

$ws1.Range("A18").Value = "04:53:50"
$ws1.Range("B18").Value = "05:15"
$ws1.Range("C18").Value = "14_ABASTO"
$ws1.Range("D18").Value = 22
$ws1.Range("E18").Value = "LP1912"

# Append two new rows at the end (rows 36 and 37)
$ws1.Range("A36").Value = "04:53:50"
$ws1.Range("B36").Value = "06:44"
$ws1.Range("C36").Value = "225_C ROCA-H SUR"
$ws1.Range("D36").Value = 111
$ws1.Range("E36").Value = "LP1912"

$ws1.Range("A37").Value = "04:53:50"
$ws1.Range("B37").Value = "06:46"
$ws1.Range("C37").Value = "215C_EL PATO"
$ws1.Range("D37").Value = 113
$ws1.Range("E37").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

# Header updates
$ws2.Range("A2").Value = "Última actualización: 04:53:50"
$ws2.Range("A3").Value = "Total filas: 10"

# Append a new row at the end (row 15)
$ws2.Range("A15").Value = "04:53:50"
$ws2.Range("B15").Value = "06:46"
$ws2.Range("C15").Value = "215C_EL PATO"
$ws2.Range("D15").Value = 113
$ws2.Range("E15").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

# Header updates
$ws3.Range("A2").Value = "Última actualización: 04:53:50"
$ws3.Range("A3").Value = "Total filas: 7"

# Insert a new row before current row 9 (old row 9 shifts to row 10)
$ws3.Range("A9").EntireRow.Insert()

$ws3.Range("A9").Value = "04:53:50"
$ws3.Range("B9").Value = "06:08"
$ws3.Range("C9").Value = "215A_LA PLATA"
$ws3.Range("D9").Value = 75
$ws3.Range("E9").Value = "L6173"

# Insert a new row before current row 11 (old row 10, now at row 10, shifts to row 12)
$ws3.Range("A11").EntireRow.Insert()

$ws3.Range("A11").Value = "04:53:50"
$ws3.Range("B11").Value = "06:32"
$ws3.Range("C11").Value = "215C_LA PLATA"
$ws3.Range("D11").Value = 99
$ws3.Range("E11").Value = "L6203"
